# Updates cryptos list values (price + volume%) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to Text so numeric-looking strings ("226.68") are not
    # auto-coerced to numbers, then restore the default (unstyled) cell style
    # so no stray NumberFormat is left behind on the cell.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.607.33"
$ws.Range("E2").Value = "  +1.22%  "

Set-TextValue "D3" "1.794.04"
$ws.Range("E3").Value = "  +0.59%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "226.68"
$ws.Range("E5").Value = "  +0.33%  "

Set-TextValue "D6" "0.558"
$ws.Range("E6").Value = "  +1.92%  "

$ws.Range("E7").Value = "  -0.02%  "

Set-TextValue "D8" "32.84"
$ws.Range("E8").Value = "  +3.27%  "

Set-TextValue "D9" "0.297"
$ws.Range("E9").Value = "  +1.86%  "

$ws.Range("E10").Value = "  +0.92%  "

Set-TextValue "D11" "0.0951"
$ws.Range("E11").Value = "  +0.44%  "

Set-TextValue "D12" "2.052.37"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("E13").Value = "  +0.83%  "

Set-TextValue "D14" "1.791.79"
$ws.Range("E14").Value = "  +0.13%  "

Set-TextValue "D15" "0.636"
$ws.Range("E15").Value = "  +2.10%  "

Set-TextValue "D16" "34.553.48"
$ws.Range("E16").Value = "  +1.33%  "

$ws.Range("E17").Value = "  +2.52%  "

Set-TextValue "D18" "68.79"
$ws.Range("E18").Value = "  +1.12%  "

Set-TextValue "D19" "248.06"
$ws.Range("E19").Value = "  +0.74%  "

Set-TextValue "D20" "0.0₃0800"
$ws.Range("E20").Value = "  +2.58%  "

Set-TextValue "D21" "11.28"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("E22").Value = "  +0.04%  "

Set-TextValue "D23" "4.18"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("E24").Value = "  +1.55%  "

Set-TextValue "D25" "166.12"
$ws.Range("E25").Value = "  +2.71%  "

Set-TextValue "D26" "7.28"
$ws.Range("E26").Value = "  +1.32%  "

Set-TextValue "D27" "16.56"
$ws.Range("E27").Value = "  +1.46%  "

Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.12%  "

Set-TextValue "D30" "4.14"
$ws.Range("E30").Value = "  +13.36%  "

$ws.Range("E31").Value = "  +2.71%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.24"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.0523"
$ws.Range("E33").Value = "  +0.54%  "

Set-TextValue "D34" "1.84"
$ws.Range("E34").Value = "  +2.21%  "

Set-TextValue "D35" "1.427.15"
$ws.Range("E35").Value = "  -1.18%  "

Set-TextValue "D36" "2.58"
$ws.Range("E36").Value = "  +5.86%  "

Set-TextValue "D37" "0.672"
$ws.Range("E37").Value = "  +2.50%  "

$ws.Range("E38").Value = "  +1.93%  "

$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("E40").Value = "  +6.34%  "

$ws.Range("E41").Value = "  +1.02%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.934"
$ws.Range("E42").Value = "  +0.99%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.75"
$ws.Range("E43").Value = "  +2.72%  "

Set-TextValue "D44" "13.67"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("E45").Value = "  +3.77%  "

$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("E47").Value = "  +0.54%  "

Set-TextValue "D48" "1.952.60"
$ws.Range("E48").Value = "  +0.44%  "

Set-TextValue "D49" "106.05"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("E51").Value = "  -5.37%  "
